$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) E21: add rewards_earned / rewards_spent / rewards_balance lines
#    to the end-user login response JSON sample.
# ------------------------------------------------------------------
$e21 = @'
{
  "loyalty_end_user_login_rs": {
    "user_info": {
      "status": "success",
      "message": "Login successful",
      "user_info": {
        "user_id": "101",
        "email": "user@example.com",
        "tier": {
          "tier_id": "3",
          "tier_name": "Gold"
        },
        "assigned_offers": [
          {
            "offer_id": "1",
            "offer_name": "Exclusive Access",
            "offer_desc": "Exclusive product launch preview"
          },
          {
            "offer_id": "2",
            "offer_name": "Priority Support",
            "offer_desc": "Enjoy priority access to our customer support team"
          }
        ],
        "wallet_info": {
          "ada_balance": "1200",
          "rewards_earned":"800",
         "rewards_spent:"300",
          "rewards_balance":"500",
          "transactions": [
            {
              "transaction_id": "1",
              "date": "2024-01-10",
              "amount": "100",
              "type": "credit",
              "desc": "ADA reward for Gold-tier spending"
            },
            {
              "transaction_id": "2",
              "date": "2024-01-15",
              "amount": "50",
              "type": "debit",
              "desc": "Purchase of product"
            }
          ]
        }
      }
    }
  }
}
'@
$ws.Range("E21").Value = $e21

# The extra lines would otherwise nudge row 21 past Excel's row-height
# cap; keep it pinned at the same capped height it already had.
$ws.Rows.Item(21).RowHeight = 409.6

# ------------------------------------------------------------------
# 2) Row 25: "create a wallet" story becomes "integrate existing wallet"
# ------------------------------------------------------------------
$ws.Range("A25").Value = "After successful login, the end user needs to integrate existing wallet"

$d25 = @'
{
  "loyalty_end_user_wallet_integrate_rq": {
    "header": {
      "user_name": "endUser",
      "product": "lrs",
      "request_type": "END_USER_CREATE_WALLET"
    },
    "wallet_info": {
      "user_id": "1",
      "currency_type": "ada",
      "wallet_name": "cardanoWallet",
     "wallet_address": "addr_test1vppvktxxw8eyhwkdf1jzq5xxqpxk8sj9d7pzvntfkng94ycn9mjxq",
    }
  }
}
'@
$ws.Range("D25").Value = $d25

$e25 = @'
{
  "loyalty_end_user_wallet_integrate_rs": {
     "status": "success"
  }
}
'@
$ws.Range("E25").Value = $e25

# Former F25 comment about the "data" field no longer applies - clear it.
$ws.Range("F25").ClearContents()

# Row 25 grew a little taller to fit the updated request sample.
$ws.Rows.Item(25).RowHeight = 245.25

# ------------------------------------------------------------------
# 3) Sheet selection moves from E9 to D9
# ------------------------------------------------------------------
$ws.Range("D9").Select()
